# Insert a new data row at row 55 (pushing the existing rows 55-77 down to
# 56-78) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55..77 down to 56..78, leaving a blank row 55 to fill in.
$ws.Rows.Item(55).Insert()

$ws.Cells.Item(55, 1).Value  = 1
$ws.Cells.Item(55, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value  = 44809
$ws.Cells.Item(55, 5).Value  = 15
$ws.Cells.Item(55, 6).Value  = 100112012
$ws.Cells.Item(55, 7).Value  = "Espinaca"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 250
$ws.Cells.Item(55, 11).Value = 3000
$ws.Cells.Item(55, 12).Value = 3500
$ws.Cells.Item(55, 13).Value = 3250
$ws.Cells.Item(55, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 1083
$ws.Cells.Item(55, 17).Value = 3
$ws.Cells.Item(55, 18).Value = "Hortaliza"
